# Scheduled-runner update: refresh Spriggan Profits price/profit figures
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets.
# Updates currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ columns (H-N) for the
# specific leve rows whose market prices changed, including a few rows where a
# column gained or lost a cached value (handled via ClearContents / new Value).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1316.3334
$ws.Range("J19").Value = 1212.25
$ws.Range("L19").Value = 1212.25
$ws.Range("N19").Value = -1562.25

$ws.Range("H40").Value = 3460.0667
$ws.Range("I40").Value = 3309.182
$ws.Range("J40").Value = 3875
$ws.Range("K40").Value = 3309.182
$ws.Range("L40").Value = 3875
$ws.Range("M40").Value = -3134.182
$ws.Range("N40").Value = -4225

$ws.Range("H113").Value = 2057.7144
$ws.Range("I113").Value = 2225.6667
$ws.Range("J113").Value = 1050
$ws.Range("K113").Value = 2225.6667
$ws.Range("L113").Value = 1050
$ws.Range("M113").Value = 1028.3333
$ws.Range("N113").Value = -7558

$ws.Range("H125").Value = 2537.0625
$ws.Range("J125").Value = 2399.4
$ws.Range("L125").Value = 21594.6
$ws.Range("N125").Value = -26514.6

$ws.Range("H137").Value = 1737
$ws.Range("I137").Value = 1345.6923
$ws.Range("K137").Value = 4037.0769
$ws.Range("M137").Value = -1487.0769

$ws.Range("H141").Value = 713.82355
$ws.Range("I141").Value = 695.96875
$ws.Range("K141").Value = 2087.90625
$ws.Range("M141").Value = 3092.09375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H69").Value = 285000
$ws.Range("J69").Value = 285000
$ws.Range("L69").Value = 285000
$ws.Range("N69").Value = -286498

$ws.Range("H72").Value = 285000
$ws.Range("J72").Value = 285000
$ws.Range("L72").Value = 855000
$ws.Range("N72").Value = -862488

$ws.Range("H102").Value = 214877.73
$ws.Range("I102").Value = 246744.4
$ws.Range("K102").Value = 246744.4
$ws.Range("M102").Value = -245122.4

$ws.Range("H132").Value = 3033523.2
$ws.Range("I132").Value = 4002623.5
$ws.Range("K132").Value = 12007870.5
$ws.Range("M132").Value = -12005340.5

$ws.Range("H139").Value = 153750
$ws.Range("J139").Value = 153750
$ws.Range("L139").Value = 153750
$ws.Range("N139").Value = -164030

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5929.1143
$ws.Range("I31").Value = 2746.9048
$ws.Range("J31").Value = 10702.429
$ws.Range("K31").Value = 2746.9048
$ws.Range("L31").Value = 10702.429
$ws.Range("M31").Value = -2451.9048
$ws.Range("N31").Value = -11292.429

$ws.Range("H34").Value = 5929.1143
$ws.Range("I34").Value = 2746.9048
$ws.Range("J34").Value = 10702.429
$ws.Range("K34").Value = 2746.9048
$ws.Range("L34").Value = 10702.429
$ws.Range("M34").Value = -2544.9048
$ws.Range("N34").Value = -11106.429

$ws.Range("H80").Value = 66666
$ws.Range("J80").Value = 66666
$ws.Range("L80").Value = 66666
$ws.Range("N80").Value = -68912

$ws.Range("H83").Value = 66666
$ws.Range("J83").Value = 66666
$ws.Range("L83").Value = 199998
$ws.Range("N83").Value = -211230

$ws.Range("H97").Value = 38995.332
$ws.Range("J97").Value = 38995.332
$ws.Range("L97").Value = 38995.332
$ws.Range("N97").Value = -40977.332

$ws.Range("H132").Value = 25001896
$ws.Range("I132").Value = 28572820
$ws.Range("K132").Value = 85718460
$ws.Range("M132").Value = -85715930

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1117.8235
$ws.Range("J107").Value = 1427.3636
$ws.Range("L107").Value = 4282.0908
$ws.Range("N107").Value = -8122.0908

$ws.Range("H140").Value = 1731.0344
$ws.Range("I140").Value = 977
$ws.Range("K140").Value = 2931
$ws.Range("M140").Value = 2249

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 12751.125
$ws.Range("I102").Value = 9002.5
$ws.Range("K102").Value = 9002.5
$ws.Range("M102").Value = -7380.5

$ws.Range("H122").Value = 7103.1665
$ws.Range("I122").Value = 5296.273
$ws.Range("J122").Value = 9942.571
$ws.Range("K122").Value = 15888.819
$ws.Range("L122").Value = 29827.713
$ws.Range("M122").Value = -13438.819
$ws.Range("N122").Value = -34727.713

$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5068.1113
$ws.Range("I7").Value = 4123.2856
$ws.Range("K7").Value = 4123.2856
$ws.Range("M7").Value = -4011.2856

$ws.Range("H22").Value = 2158.375
$ws.Range("I22").Value = 2135.7144
$ws.Range("J22").Value = 2317
$ws.Range("K22").Value = 2135.7144
$ws.Range("L22").Value = 2317
$ws.Range("M22").Value = -1840.7144
$ws.Range("N22").Value = -2907

$ws.Range("H27").Value = 2158.375
$ws.Range("I27").Value = 2135.7144
$ws.Range("J27").Value = 2317
$ws.Range("K27").Value = 2135.7144
$ws.Range("L27").Value = 2317
$ws.Range("M27").Value = -2028.7144
$ws.Range("N27").Value = -2531

$ws.Range("H34").Value = 20
$ws.Range("I34").Value = 20
$ws.Range("K34").Value = 20
$ws.Range("M34").Value = 152

$ws.Range("H61").Value = 2968.0527
$ws.Range("J61").Value = 4481.6665
$ws.Range("L61").Value = 4481.6665
$ws.Range("N61").Value = -4885.6665

$ws.Range("H113").Value = 2968.0527
$ws.Range("J113").Value = 4481.6665
$ws.Range("L113").Value = 4481.6665
$ws.Range("N113").Value = -8821.666499999999

$ws.Range("H122").Value = 4328
$ws.Range("I122").Value = 3946.5
$ws.Range("K122").Value = 11839.5
$ws.Range("M122").Value = -9389.5

$ws.Range("H126").Value = 5068.1113
$ws.Range("I126").Value = 4123.2856
$ws.Range("K126").Value = 12369.8568
$ws.Range("M126").Value = -9899.856800000001

$ws.Range("H134").Value = 20000
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 20000
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 20000
$ws.Range("M134").ClearContents()
$ws.Range("N134").Value = -30140

$ws.Range("H136").Value = 1639.75
$ws.Range("I136").Value = 1440.8462
$ws.Range("J136").Value = 2501.6667
$ws.Range("K136").Value = 4322.5386
$ws.Range("L136").Value = 7505.000100000001
$ws.Range("M136").Value = -1772.5386
$ws.Range("N136").Value = -12605.0001

$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("M137").ClearContents()
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 100025000
$ws.Range("J4").Value = 300000000
$ws.Range("L4").Value = 300000000
$ws.Range("N4").Value = -300000226

$ws.Range("H126").Value = 1191.4286
$ws.Range("I126").Value = 973.3333
$ws.Range("K126").Value = 2919.9999
$ws.Range("M126").Value = -449.9998999999998

$ws.Range("H132").Value = 7144605.5
$ws.Range("I132").Value = 8930250
$ws.Range("K132").Value = 26790750
$ws.Range("M132").Value = -26788220

$ws.Range("H133").Value = 89999
$ws.Range("J133").Value = 89999
$ws.Range("L133").Value = 89999
$ws.Range("N133").Value = -100119

$ws.Range("H136").Value = 10640545
$ws.Range("I136").Value = 11365787
$ws.Range("K136").Value = 34097361
$ws.Range("M136").Value = -34094811
